# Atualizado por script em 02-12-2023 14:46
#
# This script reproduces a scraper re-run against the Malta Premier League
# 2023-2024 betting-odds sheet:
#   - A handful of already-recorded fixtures had their row order corrected
#     (the match data in columns F:V was swapped between two adjacent rows,
#     while the leading Indice/pais/torneio/temporada/data_partida columns
#     A:E stayed put).
#   - Two brand-new fixtures were appended at the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Swap the match-data columns (F:V) between the row pairs that were
#    re-ordered by the scraper.
# ---------------------------------------------------------------------
$cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")
$swapPairs = @(@(10,11), @(19,20), @(30,31), @(46,47))

foreach ($pair in $swapPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    foreach ($col in $cols) {
        $c1 = $col + $r1
        $c2 = $col + $r2
        $tmp = $ws.Range($c1).Value2
        $ws.Range($c1).Value2 = $ws.Range($c2).Value2
        $ws.Range($c2).Value2 = $tmp
    }
}

# ---------------------------------------------------------------------
# 2) Append the two new fixtures as rows 57 and 58.
#    Copy formats from the prior last row (56) for the styled columns
#    (A: bold/bordered index, E: date-time number format) so the new
#    rows visually match the rest of the table, then fill in values.
# ---------------------------------------------------------------------
$newRows = @(
    @{
        Row = 57
        Indice = 56
        Data = 45262.58333333334
        Home = "Birkirkara"
        HomeGoals = 3
        Away = "Hamrun"
        AwayGoals = 0
        HomeOpenOdds = 3.86
        HomeOpenDate = "01/12/2023 02:15"
        HomeCloseOdds = 4.01
        HomeCloseDate = "02/12/2023 12:03"
        DrawOpenOdds = 3.01
        DrawOpenDate = "01/12/2023 02:15"
        DrawCloseOdds = 3.16
        DrawCloseDate = "02/12/2023 12:03"
        AwayOpenOdds = 1.94
        AwayOpenDate = "01/12/2023 02:15"
        AwayCloseOdds = 1.97
        AwayCloseDate = "02/12/2023 12:03"
        Url = "https://www.betexplorer.com/football/malta/premier-league/birkirkara-hamrun/dSNBZtPR/"
    },
    @{
        Row = 58
        Indice = 57
        Data = 45262.58333333334
        Home = "Santa Lucia"
        HomeGoals = 1
        Away = "Mosta"
        AwayGoals = 0
        HomeOpenOdds = 3.84
        HomeOpenDate = "01/12/2023 02:15"
        HomeCloseOdds = 2.98
        HomeCloseDate = "02/12/2023 13:57"
        DrawOpenOdds = 3.43
        DrawOpenDate = "01/12/2023 02:15"
        DrawCloseOdds = 3.4
        DrawCloseDate = "02/12/2023 13:55"
        AwayOpenOdds = 1.8
        AwayOpenDate = "01/12/2023 02:15"
        AwayCloseOdds = 2.26
        AwayCloseDate = "02/12/2023 13:57"
        Url = "https://www.betexplorer.com/football/malta/premier-league/santa-lucia-mosta-fc/25O7zN9L/"
    }
)

foreach ($nr in $newRows) {
    $r = $nr.Row

    # Clone the styled formatting of the fixed "Indice" (A) and
    # "data_partida" (E) columns from the last existing row.
    $ws.Range("A56").Copy()
    $ws.Range("A" + $r).PasteSpecial(-4122)

    $ws.Range("E56").Copy()
    $ws.Range("E" + $r).PasteSpecial(-4122)

    $ws.Range("A" + $r).Value2 = $nr.Indice
    $ws.Range("B" + $r).Value2 = "malta"
    $ws.Range("C" + $r).Value2 = "premier-league"
    $ws.Range("D" + $r).Value2 = "2023-2024"
    $ws.Range("E" + $r).Value2 = $nr.Data
    $ws.Range("F" + $r).Value2 = $nr.Home
    $ws.Range("G" + $r).Value2 = $nr.HomeGoals
    $ws.Range("H" + $r).Value2 = $nr.Away
    $ws.Range("I" + $r).Value2 = $nr.AwayGoals
    $ws.Range("J" + $r).Value2 = $nr.HomeOpenOdds
    $ws.Range("K" + $r).Value2 = $nr.HomeOpenDate
    $ws.Range("L" + $r).Value2 = $nr.HomeCloseOdds
    $ws.Range("M" + $r).Value2 = $nr.HomeCloseDate
    $ws.Range("N" + $r).Value2 = $nr.DrawOpenOdds
    $ws.Range("O" + $r).Value2 = $nr.DrawOpenDate
    $ws.Range("P" + $r).Value2 = $nr.DrawCloseOdds
    $ws.Range("Q" + $r).Value2 = $nr.DrawCloseDate
    $ws.Range("R" + $r).Value2 = $nr.AwayOpenOdds
    $ws.Range("S" + $r).Value2 = $nr.AwayOpenDate
    $ws.Range("T" + $r).Value2 = $nr.AwayCloseOdds
    $ws.Range("U" + $r).Value2 = $nr.AwayCloseDate
    $ws.Range("V" + $r).Value2 = $nr.Url
}
